$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.182.21"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "2.092.94"
$ws.Range("E3").Value = "  +9.54%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'251.83"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E6").Value = "  -3.97%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'50.66"
$ws.Range("E8").Value = "  +7.03%  "
$ws.Range("D9").Value = "'61.04"
$ws.Range("E9").Value = "  +5.09%  "
$ws.Range("D10").Value = "'0.373"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'0.0747"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("E12").Value = "  +5.98%  "
$ws.Range("D13").Value = "'15.13"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "2.394.56"
$ws.Range("E14").Value = "  +9.27%  "
$ws.Range("D15").Value = "'0.837"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").Value = "2.091.66"
$ws.Range("E16").Value = "  +9.49%  "
$ws.Range("D17").Value = "'5.13"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "37.019.44"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'72.60"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Value = "0.0₃0824"
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").Value = "'13.38"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").Value = "'240.79"
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'2.50"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "'169.94"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'21.39"
$ws.Range("E27").Value = "  +15.10%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'9.34"
$ws.Range("E28").Value = "  +6.52%  "
$ws.Range("E29").Value = "  -3.63%  "
$ws.Range("D30").Value = "'26.33"
$ws.Range("E30").Value = "  +38.67%  "
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "  +27.13%  "
$ws.Range("D33").Value = "'4.52"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").Value = "'0.0927"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'2.30"
$ws.Range("E37").Value = "  +18.30%  "
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("E39").Value = "  -3.73%  "
$ws.Range("E40").Value = "  -7.85%  "
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").Value = "'17.59"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("E43").Value = "  +6.48%  "
$ws.Range("D44").Value = "'98.02"
$ws.Range("E44").Value = "  -6.98%  "
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("D46").Value = "'0.0866"
$ws.Range("E46").Value = "  +3.95%  "
$ws.Range("D47").Value = "1.321.28"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  +6.12%  "
$ws.Range("D49").Value = "'6.92"
$ws.Range("E49").Value = "  +8.82%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.275.21"
$ws.Range("E50").Value = "  +7.97%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'2.27"
$ws.Range("E51").Value = "  -4.51%  "
